# Commit: "#8 A last idea for futur work added in report"
#
# Splits the "Future Work" bullet that currently reads:
#   "Create a profil for each person who talk (voice recognition by person)"
# into two separate paragraphs:
#   1) "Generate a voice to communicate with SAI"      <- brand-new idea, takes
#      over the original paragraph
#   2) "Create a profil for each person who talk (voice recognition by person)"
#      <- the former sentence, now living in its own (new) paragraph, together
#      with the _GoBack bookmark that used to sit right before the final ")".

$d = $word.ActiveDocument

$oldText  = "Create a profil for each person who talk (voice recognition by person)"
$newText1 = "Generate a voice to communicate with SAI"
$newText2 = "Create a profil for each person who talk (voice recognition by person)"

# Locate the paragraph holding the original sentence. Paragraph Range.Text ends
# with a trailing paragraph-mark ("`r"), so trim it before comparing.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r") -eq $oldText) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start
    $paraEnd   = $target.Range.End

    # Visible text occupies [paraStart, paraEnd - 1); the very last visible
    # character is the closing ")" which lives in its own run (together with
    # the _GoBack bookmark right before it). Splitting just before that
    # character keeps the bookmark and its run intact.
    $splitPoint = $paraEnd - 2

    # Step 1: rewrite the bulk of the paragraph (everything except the final
    # ")") into the brand-new sentence.
    $firstRunRange = $d.Range($paraStart, $splitPoint)
    $firstRunRange.Text = $newText1

    # Step 2: break the paragraph right after the new sentence. This leaves
    # the original paragraph mark (and its rsid attributes) with the first
    # paragraph, while the bookmark + remaining ")" run move into a fresh
    # paragraph.
    $breakPos = $paraStart + $newText1.Length
    $breakRange = $d.Range($breakPos, $breakPos)
    $breakRange.InsertParagraphAfter()

    # Step 3: the new paragraph now holds just the bookmark followed by a
    # single run containing ")". Replace that run's text with the full
    # original sentence so the bookmark ends up right before it.
    $secondPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text.TrimEnd("`r") -eq ")") {
            $secondPara = $para
            break
        }
    }

    if ($secondPara -ne $null) {
        $secStart = $secondPara.Range.Start
        $secEnd   = $secondPara.Range.End - 1
        $secRange = $d.Range($secStart, $secEnd)
        $secRange.Text = $newText2
    }
}
